$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.428.08'
$ws.Range("E2").Value = '  +1.36%  '
$ws.Range("D3").Value = '2.288.80'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '156.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +15,482.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '96.97'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.58%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.533'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  +1.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '35.46'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0810'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("D15").Value = '2.641.65'
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("E16").Value = '  +1.33%  '
$ws.Range("D17").Value = '2.277.16'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("E18").Value = '  +3.97%  '
$ws.Range("D19").Value = '42.273.96'
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.51%  '
$ws.Range("D21").Value = '0.0₃0920'
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '244.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.34%  '
$ws.Range("E30").Value = '  +0.89%  '
$ws.Range("E31").Value = '  +1.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '161.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0754'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("E36").Value = '  +1.60%  '
$ws.Range("E37").Value = '  +4.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.66%  '
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.116'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("E42").Value = '  +7.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.80'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.75%  '
$ws.Range("D44").Value = '2.012.82'
$ws.Range("E44").Value = '  -2.97%  '
$ws.Range("E45").Value = '  +11.13%  '
$ws.Range("E46").Value = '  +2.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("E48").Value = '  +2.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.15'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.27%  '
$ws.Range("E50").Value = '  +0.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.99'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.12%  '
